$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A24").Value = 44433
$ws.Range("B24").Value = 1
$ws.Range("D24").Value = "Research on UI/UX for chart navigation"

$ws.Range("A25").Value = 44442
$ws.Range("B25").Value = 4
$ws.Range("D25").Value = "Further implemented chart navigation"

$ws.Range("A26").Value = 44443
$ws.Range("B26").Value = 3
$ws.Range("D26").Value = "Further implemented chart navigation"

$ws.Range("L7").Select()
